$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "Correct price"

# Add rows 6-14: column A = row number, column B = "David <n>"
for ($i = 6; $i -le 14; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
    $ws.Cells.Item($i, 2).Value = "David $i"
}
